$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header H1's formatting (bold, border, centered) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-14
$values = @(
    @(5, 6),
    @(6, 7),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
